# Rename the three header/footer logo pictures (swap their generated
# "imageN.ext" display names), matching the authored XML diff:
#   - header  -> BTEC logo:    image2.jpg -> image1.jpg
#   - footer1 -> Pearson logo: image1.png -> image2.png  (id=3)
#   - footer2 -> Pearson logo: image1.png -> image2.png  (id=2)
#
# These are inline pictures living in the section's headers/footers, so we
# reach them through Sections(1).Headers/Footers rather than
# ActiveDocument.InlineShapes (the body has none).

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# BTec_Logo-Orange picture, in the header that carries content (Headers(2)).
$btecHeader = $sec.Headers.Item(2)
$btecShape = $btecHeader.Range.InlineShapes.Item(1)
$btecShape.Name = "image1.jpg"

# PearsonLogo picture in the footer that carries content (Footers(2)).
$pearsonFooterA = $sec.Footers.Item(2)
$pearsonShapeA = $pearsonFooterA.Range.InlineShapes.Item(1)
$pearsonShapeA.Name = "image2.png"

# PearsonLogo picture in the other content-bearing footer (Footers(1)).
$pearsonFooterB = $sec.Footers.Item(1)
$pearsonShapeB = $pearsonFooterB.Range.InlineShapes.Item(1)
$pearsonShapeB.Name = "image2.png"

Write-Output ("BTec shape name: " + $btecShape.Name)
Write-Output ("Pearson footer2 shape name: " + $pearsonShapeA.Name)
Write-Output ("Pearson footer1 shape name: " + $pearsonShapeB.Name)
